# feat: add 2022-Q4 data
#
# 1. Insert a brand-new "2022-Q4" worksheet right after "总计" (so it sits
#    before "2022-Q3" in tab order), cloning the layout/formatting of the
#    existing "2022-Q3" sheet and filling in the new quarter's numbers.
# 2. Insert a matching summary row into "总计" for the new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q4" sheet, positioned after "总计".
# ---------------------------------------------------------------------
$zongji = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $zongji)
$newSheet.Name = "2022-Q4"

# Clone structure/formatting/values from "2022-Q3" (the engine re-resolves
# worksheet handles live by position/name, so fetch this fresh, after the
# Add() call, rather than reusing a handle obtained earlier).
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Range("A1:H8").Copy($wb.Worksheets.Item("2022-Q4").Range("A1"))

# ---------------------------------------------------------------------
# Step 2: overwrite the new sheet's per-fund numbers with the Q4 figures.
# Columns D/E/F/G are stored as text (like every other quarter sheet), so
# assign them with a leading apostrophe to keep them as text instead of
# letting Excel auto-coerce the numeric-looking strings into numbers.
# Column H (ranking) stays a genuine number.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2022-Q4")

$q4Data = @(
    @{ Row = 2; D = "3.27"; E = "90.95"; F = "7.01"; G = "0.2292"; H = 4 },
    @{ Row = 3; D = "3.25"; E = "94.43"; F = "6.65"; G = "0.2161"; H = 4 },
    @{ Row = 4; D = "3.06"; E = "94.24"; F = "5.64"; G = "0.1726"; H = 5 },
    @{ Row = 5; D = "3.23"; E = "83.79"; F = "5.34"; G = "0.1725"; H = 7 },
    @{ Row = 6; D = "0.81"; E = "83.79"; F = "5.34"; G = "0.0433"; H = 7 },
    @{ Row = 7; D = "0.45"; E = "94.43"; F = "6.65"; G = "0.0299"; H = 4 },
    @{ Row = 8; D = "0.48"; E = "94.24"; F = "5.64"; G = "0.0271"; H = 5 }
)

foreach ($entry in $q4Data) {
    $r = $entry.Row
    $q4.Range("D$r").Value2 = "'" + $entry.D
    $q4.Range("E$r").Value2 = "'" + $entry.E
    $q4.Range("F$r").Value2 = "'" + $entry.F
    $q4.Range("G$r").Value2 = "'" + $entry.G
    $q4.Range("H$r").Value2 = $entry.H
}

# ---------------------------------------------------------------------
# Step 3: add the new "2022-Q4" row to the "总计" (totals) summary sheet.
# Shift the existing Q3/Q2/Q1 rows down one, then write the new row 2.
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

$zj.Range("A4:D4").Copy($zj.Range("A5"))
$zj.Range("A3:D3").Copy($zj.Range("A4"))
$zj.Range("A2:D2").Copy($zj.Range("A3"))

$zj.Range("A2").Value2 = 0
$zj.Range("B2").Value2 = "2022-Q4"
$zj.Range("C2").Value2 = 7
$zj.Range("D2").Value2 = 0.89

$zj.Range("A3").Value2 = 1
$zj.Range("A4").Value2 = 2
$zj.Range("A5").Value2 = 3
